$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E31: "n/a" -> "check" (chapter 4 index term status)
$ws.Range("E31").Value = "check"

# Add Chapter 4 (column E) status for rows 64-77
$ws.Range("E64").Value = "check"
$ws.Range("E65").Value = "check"
$ws.Range("E66").Value = "check"
$ws.Range("E67").Value = "check"
$ws.Range("E68").Value = "check"
$ws.Range("E69").Value = "check"
$ws.Range("E70").Value = "n/a"
$ws.Range("E71").Value = "check"
$ws.Range("E72").Value = "check"
$ws.Range("E73").Value = "n/a"
$ws.Range("E74").Value = "n/a"
$ws.Range("E75").Value = "n/a"
$ws.Range("E76").Value = "n/a"
$ws.Range("E77").Value = "check"

# Adjust page layout / view: scroll back to top-left and move selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E76").Select()
